$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "Which of the following is not considered to be a pointing device?",
        "ques_type": 2,
        "options": [
            "Keyboard",
            "Touchpad",
            "Stylus",
            "Mouse"
        ],
        "score": "Keyboard"
    },
    {
        "title": "Which of the following is a browser used to access the internet?",
        "ques_type": 2,
        "options": [
            "Google Hangouts",
            "Google Chrome",
            "Google Search",
            "Google Plus"
        ],
        "score": "Google Chrome"
    },
    {
        "title": "Which of the following is recommended to ensure email safety?",
        "ques_type": 2,
        "options": [
            "Clicking on attachments from unknown senders.",
            "Never logging out of your account.",
            "Uninstalling spam filters.",
            "Avoiding public Wi-Fi."
        ],
        "score": "Avoiding public Wi-Fi."
    },
    {
        "title": "You oversee numerous projects with corresponding files on your desktop, leading to clutter and difficulty in locating specific documents. Your goal is to organize these files to enhance efficiency and traceability.What is the best method to organize your desktop files?",
        "ques_type": 2,
        "options": [
            "Create project-specific folders and subfolders based on document types.",
            "Sort all files by size and date, then place them in general folders.",
            "Create a folder for each team member and distribute files accordingly. ",
            "Arrange files alphabetically in a single folder for easy access."
        ],
        "score": "Create project-specific folders and subfolders based on document types."
    }
]
'@

# Remove the old row 2 (clears value + shifts dimension back to A1)
$ws.Range("A2").EntireRow.Delete() | Out-Null

# Clear any special formatting (border, bold, alignment) on A1 so it reverts to default style
$ws.Range("A1").ClearFormats() | Out-Null

# Put the new JSON-like text into A1
$ws.Range("A1").Value = $newText
$ws.Rows.Item(1).AutoFit() | Out-Null

